$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new blank rows right after the header row (row 1), shifting the
# existing data (previously rows 2-32) down to rows 13-43.
$ws.Rows("2:12").Insert()

# The freshly inserted rows don't carry the same cell formatting as the rest
# of the data table (date format with border on column A). Copy the
# formatting from the row that used to be row 2 (now row 13) down onto the
# newly inserted rows so they match the rest of the series.
$ws.Range("A13:B13").Copy()
$ws.Range("A2:B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "backward extension" rows with the real-time data values.
$ws.Range("A2").Value = 30681
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 31047
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 31412
$ws.Range("B4").Value = [double]"2.220446049250313E-14"

$ws.Range("A5").Value = 31777
$ws.Range("B5").Value = [double]"-2.220446049250313E-14"

$ws.Range("A6").Value = 32142
$ws.Range("B6").Value = [double]"-2.220446049250313E-14"

$ws.Range("A7").Value = 32508
$ws.Range("B7").Value = [double]"-2.220446049250313E-14"

$ws.Range("A8").Value = 32873
$ws.Range("B8").Value = [double]"2.220446049250313E-14"

$ws.Range("A9").Value = 33238
$ws.Range("B9").Value = [double]"2.220446049250313E-14"

$ws.Range("A10").Value = 33603
$ws.Range("B10").Value = 0.8650282515740848

$ws.Range("A11").Value = 33969
$ws.Range("B11").Value = 0.2387091425554155

$ws.Range("A12").Value = 34334
$ws.Range("B12").Value = -0.2106527079487774
